$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.419.31"
$ws.Range("E2").Value = "'  +1.72%  "
$ws.Range("D3").Value = "'3.569.71"
$ws.Range("E3").Value = "'  +2.27%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "'  +0.00%  "
$ws.Range("D5").Value = "'621.88"
$ws.Range("E5").Value = "'  +2.78%  "
$ws.Range("D6").Value = "'155.02"
$ws.Range("E6").Value = "'  +4.34%  "
$ws.Range("D7").Value = "'3.566.71"
$ws.Range("E7").Value = "'  +2.27%  "
$ws.Range("E8").Value = "'  -0.01%  "
$ws.Range("E9").Value = "'  +2.25%  "
$ws.Range("E10").Value = "'  +5.65%  "
$ws.Range("E11").Value = "'  +5.43%  "
$ws.Range("D12").Value = "'0.438"
$ws.Range("E12").Value = "'  +3.88%  "
$ws.Range("D13").Value = "'0.0000221"
$ws.Range("E13").Value = "'  +1.70%  "
$ws.Range("D14").Value = "'33.13"
$ws.Range("E14").Value = "'  +5.50%  "
$ws.Range("D15").Value = "'4.169.99"
$ws.Range("E15").Value = "'  +2.16%  "
$ws.Range("D16").Value = "'3.563.44"
$ws.Range("E16").Value = "'  +2.21%  "
$ws.Range("D17").Value = "'68.288.52"
$ws.Range("E17").Value = "'  +1.51%  "
$ws.Range("E18").Value = "'  +0.00%  "
$ws.Range("D19").Value = "'6.77"
$ws.Range("E19").Value = "'  +5.99%  "
$ws.Range("D20").Value = "'16.00"
$ws.Range("E20").Value = "'  +6.47%  "
$ws.Range("D21").Value = "'10.01"
$ws.Range("E21").Value = "'  +11.13%  "
$ws.Range("D22").Value = "'455.61"
$ws.Range("E22").Value = "'  +2.11%  "
$ws.Range("E23").Value = "'  +3.74%  "
$ws.Range("D24").Value = "'78.74"
$ws.Range("E24").Value = "'  +2.17%  "
$ws.Range("E25").Value = "'  +2.38%  "
$ws.Range("D26").Value = "'3.706.46"
$ws.Range("E26").Value = "'  +2.10%  "
$ws.Range("E27").Value = "'  -0.07%  "
$ws.Range("D28").Value = "'10.55"
$ws.Range("E28").Value = "'  +4.30%  "
$ws.Range("D29").Value = "'9.09"
$ws.Range("E29").Value = "'  +9.99%  "
$ws.Range("B30").Value = "'Fetch.AI"
$ws.Range("C30").Value = "'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D30").Value = "'1.70"
$ws.Range("E30").Value = "'  +9.37%  "
$ws.Range("B31").Value = "'PancakeSwap"
$ws.Range("C31").Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'2.57"
$ws.Range("E31").Value = "'  +3.90%  "
$ws.Range("E32").Value = "'  +5.82%  "
$ws.Range("D33").Value = "'0.999"
$ws.Range("E33").Value = "'  -0.08%  "
$ws.Range("D34").Value = "'6.35"
$ws.Range("E34").Value = "'  +3.60%  "
$ws.Range("D35").Value = "'26.16"
$ws.Range("E35").Value = "'  +2.13%  "
$ws.Range("E36").Value = "'  +3.68%  "
$ws.Range("D37").Value = "'3.561.11"
$ws.Range("E37").Value = "'  +2.32%  "
$ws.Range("D38").Value = "'8.28"
$ws.Range("E38").Value = "'  +3.59%  "
$ws.Range("E39").Value = "'  +8.97%  "
$ws.Range("E40").Value = "'  -0.02%  "
$ws.Range("D41").Value = "'178.42"
$ws.Range("E41").Value = "'  +2.65%  "
$ws.Range("E42").Value = "'  +4.97%  "
$ws.Range("D43").Value = "'0.999"
$ws.Range("E43").Value = "'  -0.04%  "
$ws.Range("E44").Value = "'  +3.50%  "
$ws.Range("D45").Value = "'30.94"
$ws.Range("E45").Value = "'  +14.38%  "
$ws.Range("D46").Value = "'0.899"
$ws.Range("E46").Value = "'  +1.97%  "
$ws.Range("E47").Value = "'  +2.22%  "
$ws.Range("E48").Value = "'  +6.96%  "
$ws.Range("E49").Value = "'  +4.34%  "
$ws.Range("D50").Value = "'7.80"
$ws.Range("E50").Value = "'  +3.58%  "
$ws.Range("E51").Value = "'  +7.61%  "
